$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 2
Set-TextValue "D2" "22.477.73"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
Set-TextValue "D3" "1.575.69"
$ws.Range("E3").Value = "  +1.02%  "

# Row 4
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("E5").Value = "  -0.11%  "

# Row 6
$ws.Range("E6").Value = "  +0.83%  "

# Row 7
Set-TextValue "D7" "0.3694"
$ws.Range("E7").Value = "  +1.55%  "

# Row 8
Set-TextValue "D8" "47.84"
$ws.Range("E8").Value = "  -1.70%  "

# Row 9
$ws.Range("E9").Value = "  -0.12%  "

# Row 10
$ws.Range("E10").Value = "  +2.56%  "

# Row 11
Set-TextValue "D11" "0.07567"
$ws.Range("E11").Value = "  +2.68%  "

# Row 12
$ws.Range("E12").Value = "  -0.08%  "

# Row 13
$ws.Range("E13").Value = "  +0.72%  "

# Row 14
Set-TextValue "D14" "5.949"
$ws.Range("E14").Value = "  +0.94%  "

# Row 15
Set-TextValue "D15" "6.949"
$ws.Range("E15").Value = "  +1.57%  "

# Row 16
Set-TextValue "D16" "1.568.80"
$ws.Range("E16").Value = "  +0.47%  "

# Row 17
Set-TextValue "D17" "0.00001121"
$ws.Range("E17").Value = "  +2.16%  "

# Row 18
Set-TextValue "D18" "88.30"
$ws.Range("E18").Value = "  -0.36%  "

# Row 19
Set-TextValue "D19" "0.06732"
$ws.Range("E19").Value = "  -0.04%  "

# Row 20
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D20" "0.9998"
$ws.Range("E20").Value = "  -0.14%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "6.388"
$ws.Range("E21").Value = "  +1.61%  "

# Row 22
Set-TextValue "D22" "16.53"
$ws.Range("E22").Value = "  +3.57%  "

# Row 23
Set-TextValue "D23" "12.04"
$ws.Range("E23").Value = "  +1.17%  "

# Row 24
Set-TextValue "D24" "22.465.60"
$ws.Range("E24").Value = "  +0.43%  "

# Row 25
Set-TextValue "D25" "2.386"
$ws.Range("E25").Value = "  -0.29%  "

# Row 26
Set-TextValue "D26" "2.640"
$ws.Range("E26").Value = "  +3.01%  "

# Row 27
Set-TextValue "D27" "150.94"
$ws.Range("E27").Value = "  +1.05%  "

# Row 28
$ws.Range("E28").Value = "  +1.95%  "

# Row 29
Set-TextValue "D29" "4.988"
$ws.Range("E29").Value = "  -0.22%  "

# Row 30
Set-TextValue "D30" "125.56"
$ws.Range("E30").Value = "  +2.34%  "

# Row 31
Set-TextValue "D31" "1.746.61"
$ws.Range("E31").Value = "  +0.51%  "

# Row 32
Set-TextValue "D32" "1.091"
$ws.Range("E32").Value = "  +3.27%  "

# Row 33
Set-TextValue "D33" "6.113"
$ws.Range("E33").Value = "  +0.47%  "

# Row 34
Set-TextValue "D34" "1.997"
$ws.Range("E34").Value = "  +0.39%  "

# Row 35
Set-TextValue "D35" "9.869"
$ws.Range("E35").Value = "  +3.40%  "

# Row 36
Set-TextValue "D36" "0.08382"
$ws.Range("E36").Value = "  +2.00%  "

# Row 37
Set-TextValue "D37" "0.02464"
$ws.Range("E37").Value = "  +4.09%  "

# Row 38
$ws.Range("E38").Value = "  +1.35%  "

# Row 39
Set-TextValue "D39" "0.06396"
$ws.Range("E39").Value = "  +0.85%  "

# Row 40
Set-TextValue "D40" "1.299"
$ws.Range("E40").Value = "  -0.94%  "

# Row 41
Set-TextValue "D41" "5.360"
$ws.Range("E41").Value = "  +1.02%  "

# Row 42
Set-TextValue "D42" "11.48"
$ws.Range("E42").Value = "  +3.71%  "

# Row 43
Set-TextValue "D43" "0.6279"
$ws.Range("E43").Value = "  +4.17%  "

# Row 44
Set-TextValue "D44" "14.09"
$ws.Range("E44").Value = "  +4.04%  "

# Row 45
Set-TextValue "D45" "1.000"
$ws.Range("E45").Value = "  +0.00%  "

# Row 46
Set-TextValue "D46" "0.6109"
$ws.Range("E46").Value = "  +7.03%  "

# Row 47
Set-TextValue "D47" "3.780"
$ws.Range("E47").Value = "  +0.59%  "

# Row 48
$ws.Range("E48").Value = "  +3.00%  "

# Row 49
Set-TextValue "D49" "125.37"
$ws.Range("E49").Value = "  +0.76%  "

# Row 50
Set-TextValue "D50" "1.212"
$ws.Range("E50").Value = "  +0.67%  "

# Row 51
Set-TextValue "D51" "0.07223"
